# Apply the "Add files via upload" update to the Saldo/Export sheet.
# Net effect (per the target diff):
#   - Add a new account row (004481463, Mara, 10000) before the Rafael row.
#   - Move/replace the Thomas row (008026942): delete it from its old spot
#     (value 0.34, just after the Adriano/0.35 row) and re-insert it earlier
#     in the list (before the Ahmad row) with an updated value of 2000.34.
#   - Move/replace the Camila row (004368628): delete it from its old spot
#     (value 573.28, just before the Manuela row) and re-insert it earlier
#     in the list (before the Lohran row) with an updated value of 1203.18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insertions (performed top-to-bottom so row numbers below are still valid) ---

# New row 3: 004481463 / Mara / 10000  (pushes the existing Rafael row, and
# everything after it, down by one)
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004481463"
$ws.Cells.Item(3, 2).Value = "Mara"
$ws.Cells.Item(3, 3).Value = 10000

# New row 5: 008026942 / Thomas / 2000.34 (before the Ahmad row, which is now
# at row 5 after the insertion above)
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "008026942"
$ws.Cells.Item(5, 2).Value = "Thomas"
$ws.Cells.Item(5, 3).Value = 2000.34

# New row 7: 004368628 / Camila / 1203.18 (before the Lohran row, which is now
# at row 7 after the insertions above)
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "004368628"
$ws.Cells.Item(7, 2).Value = "Camila"
$ws.Cells.Item(7, 3).Value = 1203.18

# --- Deletions of the original rows that were moved above ---

# Old Camila row (573.28), originally right before Manuela; after the three
# inserts above it is now at row 10.
$ws.Rows.Item(10).Delete()

# Old Thomas row (0.34), originally right after the Adriano/0.35 row; after
# the three inserts above (and before the single deletion above it, which
# does not affect rows below it) it is now at row 169.
$ws.Rows.Item(169).Delete()
